# Updates market-price-derived columns (currentAveragePrice / NQ / HQ,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) on several rows across all eight job
# sheets to reflect a refreshed data pull. A handful of rows also have
# their LevePriceHQ / LeveProfitHQ cells cleared entirely (no HQ price
# data available any more for those leves).

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 503.84616
$ws.Range("J2").Value = 522.7273
$ws.Range("L2").Value = 522.7273
$ws.Range("N2").Value = -748.7273
$ws.Range("H33").Value = 737.1667
$ws.Range("I33").Value = 768.25
$ws.Range("J33").Value = 302
$ws.Range("K33").Value = 768.25
$ws.Range("L33").Value = 302
$ws.Range("M33").Value = -539.25
$ws.Range("N33").Value = -760
$ws.Range("H64").Value = 3109.85
$ws.Range("I64").Value = 3099.2
$ws.Range("J64").Value = 3141.8
$ws.Range("K64").Value = 3099.2
$ws.Range("L64").Value = 3141.8
$ws.Range("M64").Value = -2851.2
$ws.Range("N64").Value = -3637.8
$ws.Range("H67").Value = 3109.85
$ws.Range("I67").Value = 3099.2
$ws.Range("J67").Value = 3141.8
$ws.Range("K67").Value = 3099.2
$ws.Range("L67").Value = 3141.8
$ws.Range("M67").Value = -2241.2
$ws.Range("N67").Value = -4857.8
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").ClearContents()
$ws.Range("N117").ClearContents()
$ws.Range("H125").Value = 1291.4546
$ws.Range("J125").Value = 1492.8889
$ws.Range("L125").Value = 13436.0001
$ws.Range("N125").Value = -18356.0001

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1156.6666
$ws.Range("I4").Value = 925.7143
$ws.Range("J4").Value = 1480
$ws.Range("K4").Value = 925.7143
$ws.Range("L4").Value = 1480
$ws.Range("M4").Value = -809.7143
$ws.Range("N4").Value = -1712
$ws.Range("H45").Value = 2544.0952
$ws.Range("I45").Value = 2663.25
$ws.Range("J45").Value = 2162.8
$ws.Range("K45").Value = 2663.25
$ws.Range("L45").Value = 2162.8
$ws.Range("M45").Value = -2286.25
$ws.Range("N45").Value = -2916.8
$ws.Range("H63").Value = 3170
$ws.Range("I63").Value = 2356.6667
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 2356.6667
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -1670.6667
$ws.Range("N63").Value = -6372
$ws.Range("H66").Value = 3170
$ws.Range("I66").Value = 2356.6667
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 11783.3335
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -8351.333500000001
$ws.Range("N66").Value = -31864
$ws.Range("H109").Value = 28725.666
$ws.Range("J109").Value = 28725.666
$ws.Range("L109").Value = 28725.666
$ws.Range("N109").Value = -31499.666
$ws.Range("H122").Value = 2767.6743
$ws.Range("I122").Value = 2454.5938
$ws.Range("J122").Value = 3678.4546
$ws.Range("K122").Value = 7363.7814
$ws.Range("L122").Value = 11035.3638
$ws.Range("M122").Value = -4913.7814
$ws.Range("N122").Value = -15935.3638

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6252232.5
$ws.Range("I86").Value = 7694485.5
$ws.Range("K86").Value = 7694485.5
$ws.Range("M86").Value = -7693362.5
$ws.Range("H89").Value = 6252232.5
$ws.Range("I89").Value = 7694485.5
$ws.Range("K89").Value = 38472427.5
$ws.Range("M89").Value = -38466811.5

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2937.6667
$ws.Range("I62").Value = 2463.125
$ws.Range("J62").Value = 3229.6924
$ws.Range("K62").Value = 2463.125
$ws.Range("L62").Value = 3229.6924
$ws.Range("M62").Value = -1839.125
$ws.Range("N62").Value = -4477.6924
$ws.Range("H65").Value = 2937.6667
$ws.Range("I65").Value = 2463.125
$ws.Range("J65").Value = 3229.6924
$ws.Range("K65").Value = 12315.625
$ws.Range("L65").Value = 16148.462
$ws.Range("M65").Value = -9195.625
$ws.Range("N65").Value = -22388.462
$ws.Range("H122").Value = 1306.65
$ws.Range("I122").Value = 891.1
$ws.Range("J122").Value = 1722.2
$ws.Range("K122").Value = 2673.3
$ws.Range("L122").Value = 5166.6
$ws.Range("M122").Value = -223.3000000000002
$ws.Range("N122").Value = -10066.6

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 848.62164
$ws.Range("J5").Value = 1234.5
$ws.Range("L5").Value = 3703.5
$ws.Range("N5").Value = -3927.5
$ws.Range("H120").Value = 10099.75
$ws.Range("I120").Value = 6999.5
$ws.Range("K120").Value = 20998.5
$ws.Range("M120").Value = -16160.5
$ws.Range("H135").Value = 848.62164
$ws.Range("J135").Value = 1234.5
$ws.Range("L135").Value = 11110.5
$ws.Range("N135").Value = -16180.5

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 82946280
$ws.Range("I70").Value = 414705900
$ws.Range("J70").Value = 6375
$ws.Range("K70").Value = 414705900
$ws.Range("L70").Value = 6375
$ws.Range("M70").Value = -414705630
$ws.Range("N70").Value = -6915
$ws.Range("H73").Value = 82946280
$ws.Range("I73").Value = 414705900
$ws.Range("J73").Value = 6375
$ws.Range("K73").Value = 414705900
$ws.Range("L73").Value = 6375
$ws.Range("M73").Value = -414704964
$ws.Range("N73").Value = -8247
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").ClearContents()
$ws.Range("N111").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("N123").ClearContents()

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2155.8333
$ws.Range("I7").Value = 1597
$ws.Range("J7").Value = 4950
$ws.Range("K7").Value = 1597
$ws.Range("L7").Value = 4950
$ws.Range("M7").Value = -1485
$ws.Range("N7").Value = -5174
$ws.Range("H122").Value = 4555.8213
$ws.Range("I122").Value = 5230.933
$ws.Range("J122").Value = 3776.8462
$ws.Range("K122").Value = 15692.799
$ws.Range("L122").Value = 11330.5386
$ws.Range("M122").Value = -13242.799
$ws.Range("N122").Value = -16230.5386
$ws.Range("H126").Value = 2155.8333
$ws.Range("I126").Value = 1597
$ws.Range("J126").Value = 4950
$ws.Range("K126").Value = 4791
$ws.Range("L126").Value = 14850
$ws.Range("M126").Value = -2321
$ws.Range("N126").Value = -19790

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 682.5294
$ws.Range("I122").Value = 682.5294
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2047.5882
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 402.4117999999999
$ws.Range("N122").ClearContents()
